$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update provider e-mail addresses (B2, B3, B4, B6) and turn them into mailto hyperlinks
$ws.Range("B2").Value = "japsequiposelectricos@gmail.com"
$ws.Range("B3").Value = "japsequiposelectricos@gmail.com"
$ws.Range("B4").Value = "japsequiposelectricos@gmail.com"
$ws.Range("B6").Value = "japsequiposelectricos@gmail.com"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:japsequiposelectricos@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:japsequiposelectricos@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:japsequiposelectricos@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:japsequiposelectricos@gmail.com")

# Remove the duplicated row 7 entry (A7/B7) but keep the formatted/empty B7 cell and G7
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()

# Move the active selection to B7
$ws.Range("B7").Select()
